# Apply odds updates to the active worksheet (Sheet1) of the Betfair
# Back/Lay workbook for 2026-01-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Tigres vs Pumas UNAM)
$ws.Range("W2").Value = 2.42
$ws.Range("AA2").Value = 220

# Row 3 (Toluca vs Santos Laguna)
$ws.Range("H3").Value = 14
$ws.Range("O3").Value = 1.1
$ws.Range("P3").Value = 3.65
$ws.Range("Q3").Value = 1.29
$ws.Range("R3").Value = 2.12
$ws.Range("S3").Value = 1.76
$ws.Range("T3").Value = 1.78
$ws.Range("W3").Value = 4.9
$ws.Range("X3").Value = 55

# Row 4 (Verona vs Bologna)
$ws.Range("F4").Value = 3.9
$ws.Range("G4").Value = 3.95
$ws.Range("I4").Value = 2.26
$ws.Range("P4").Value = 1.68

# Row 5 (Augsburg vs Union Berlin)
$ws.Range("F5").Value = 2.74
$ws.Range("G5").Value = 2.76
$ws.Range("I5").Value = 2.94
$ws.Range("R5").Value = 1.27
$ws.Range("AI5").Value = 55

# Row 6 (Como vs AC Milan)
$ws.Range("G6").Value = 3.25
$ws.Range("Z6").Value = 16.5
$ws.Range("AH6").Value = 19
$ws.Range("AI6").Value = 50
$ws.Range("AJ6").Value = 65
$ws.Range("AL6").Value = 65
